# Add new stop words (currency/finance acronyms) to the "stop_words" sheet
# and make that sheet the active tab, matching the author's edit.

$wb = $excel.ActiveWorkbook

$wsStop = $wb.Worksheets.Item("stop_words")

$newWords = @("www", "ag", "chf", "usd", "gbp", "jpy", "eur")

$startRow = 11
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $wsStop.Range("A$row").Value = $newWords[$i]
}

# Select a new "next" cell below the data, mirroring the author's workbook state
$wsStop.Range("A18").Select()

# Make stop_words the active/visible sheet when the workbook is reopened
$wsStop.Activate()

$wb.Save()
